$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 8000
$ws.Range("J56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("N56").Value = -9484
$ws.Range("H134").Value = 20000
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 20000
$ws.Range("N134").Value = -30140
$ws.Range("H135").Value = 48000
$ws.Range("J135").Value = 48000
$ws.Range("L135").Value = 48000
$ws.Range("N135").Value = -58140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5227
$ws.Range("I5").Value = 454
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 454
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -341
$ws.Range("N5").Value = -10226
$ws.Range("H7").Value = 762681.25
$ws.Range("J7").Value = 883000
$ws.Range("L7").Value = 883000
$ws.Range("N7").Value = -883226
$ws.Range("H23").Value = 637.6667
$ws.Range("I23").Value = 500
$ws.Range("J23").Value = 706.5
$ws.Range("K23").Value = 500
$ws.Range("L23").Value = 706.5
$ws.Range("M23").Value = -217
$ws.Range("N23").Value = -1272.5
$ws.Range("H38").Value = 9000
$ws.Range("J38").Value = 9000
$ws.Range("L38").Value = 9000
$ws.Range("N38").Value = -9832
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30994
$ws.Range("H56").Value = 8000
$ws.Range("J56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("N56").Value = -9478
$ws.Range("H62").Value = 27317.5
$ws.Range("I62").Value = 13000
$ws.Range("J62").Value = 30181
$ws.Range("K62").Value = 13000
$ws.Range("L62").Value = 30181
$ws.Range("M62").Value = -12314
$ws.Range("N62").Value = -31553
$ws.Range("H65").Value = 27317.5
$ws.Range("I65").Value = 13000
$ws.Range("J65").Value = 30181
$ws.Range("K65").Value = 39000
$ws.Range("L65").Value = 90543
$ws.Range("M65").Value = -35568
$ws.Range("N65").Value = -97407
$ws.Range("H68").Value = 36147.5
$ws.Range("J68").Value = 36147.5
$ws.Range("L68").Value = 36147.5
$ws.Range("N68").Value = -37769.5
$ws.Range("H71").Value = 36147.5
$ws.Range("J71").Value = 36147.5
$ws.Range("L71").Value = 108442.5
$ws.Range("N71").Value = -116554.5
$ws.Range("H76").Value = 8142.5
$ws.Range("I76").Value = 8142.5
$ws.Range("K76").Value = 8142.5
$ws.Range("M76").Value = -7827.5
$ws.Range("H79").Value = 8142.5
$ws.Range("I79").Value = 8142.5
$ws.Range("K79").Value = 8142.5
$ws.Range("M79").Value = -7050.5
$ws.Range("H81").Value = 12773.385
$ws.Range("J81").Value = 12773.385
$ws.Range("L81").Value = 12773.385
$ws.Range("N81").Value = -14895.385
$ws.Range("H82").Value = 25854.643
$ws.Range("I82").Value = 13780
$ws.Range("J82").Value = 27867.084
$ws.Range("K82").Value = 13780
$ws.Range("L82").Value = 27867.084
$ws.Range("M82").Value = -13397
$ws.Range("N82").Value = -28633.084
$ws.Range("H84").Value = 12773.385
$ws.Range("J84").Value = 12773.385
$ws.Range("L84").Value = 38320.155
$ws.Range("N84").Value = -48928.155
$ws.Range("H85").Value = 25854.643
$ws.Range("I85").Value = 13780
$ws.Range("J85").Value = 27867.084
$ws.Range("K85").Value = 13780
$ws.Range("L85").Value = 27867.084
$ws.Range("M85").Value = -12454
$ws.Range("N85").Value = -30519.084
$ws.Range("H92").Value = 20401
$ws.Range("J92").Value = 20401
$ws.Range("L92").Value = 20401
$ws.Range("N92").Value = -25393
$ws.Range("H110").Value = 28393.5
$ws.Range("J110").Value = 28393.5
$ws.Range("L110").Value = 28393.5
$ws.Range("N110").Value = -36573.5
$ws.Range("H124").Value = 59260
$ws.Range("J124").Value = 59260
$ws.Range("L124").Value = 59260
$ws.Range("N124").Value = -69080
$ws.Range("H125").Value = 52890
$ws.Range("J125").Value = 52890
$ws.Range("L125").Value = 52890
$ws.Range("N125").Value = -62730
$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880
$ws.Range("H130").Value = 274945
$ws.Range("J130").Value = 274945
$ws.Range("L130").Value = 274945
$ws.Range("N130").Value = -284985
$ws.Range("H132").Value = 31666.666
$ws.Range("J132").Value = 31666.666
$ws.Range("L132").Value = 31666.666
$ws.Range("N132").Value = -41786.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 500
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -206
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = -10776
$ws.Range("H38").Value = 2500
$ws.Range("I38").Value = 2000
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -1623
$ws.Range("N38").Value = -3754
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -10320
$ws.Range("H46").Value = 2500
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1789
$ws.Range("N46").Value = -3422
$ws.Range("H50").Value = 8742.643
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 8742.643
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 8742.643
$ws.Range("M50").Value = ""
$ws.Range("N50").Value = -9992.643
$ws.Range("H55").Value = 7950
$ws.Range("J55").Value = 7950
$ws.Range("L55").Value = 7950
$ws.Range("N55").Value = -8580
$ws.Range("H63").Value = 32978.25
$ws.Range("J63").Value = 32978.25
$ws.Range("L63").Value = 32978.25
$ws.Range("N63").Value = -34350.25
$ws.Range("H66").Value = 32978.25
$ws.Range("J66").Value = 32978.25
$ws.Range("L66").Value = 98934.75
$ws.Range("N66").Value = -105798.75
$ws.Range("H74").Value = 13255.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13255.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13255.8
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -15003.8
$ws.Range("H77").Value = 13255.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13255.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 39767.39999999999
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -48503.39999999999
$ws.Range("H132").Value = 6413541.5
$ws.Range("I132").Value = 11495204
$ws.Range("J132").Value = 6228.174
$ws.Range("K132").Value = 34485612
$ws.Range("L132").Value = 18684.522
$ws.Range("M132").Value = -34483082
$ws.Range("N132").Value = -23744.522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2821.6904
$ws.Range("I113").Value = 406.625
$ws.Range("J113").Value = 10549.9
$ws.Range("K113").Value = 1219.875
$ws.Range("L113").Value = 31649.7
$ws.Range("M113").Value = 950.125
$ws.Range("N113").Value = -35989.7
$ws.Range("H132").Value = 1994.8687
$ws.Range("I132").Value = 540.5484
$ws.Range("K132").Value = 4864.9356
$ws.Range("M132").Value = -2334.9356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 48333.332
$ws.Range("J134").Value = 48333.332
$ws.Range("L134").Value = 144999.996
$ws.Range("N134").Value = -150069.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 725
$ws.Range("I30").Value = 700
$ws.Range("K30").Value = 700
$ws.Range("M30").Value = -592
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -664
$ws.Range("N35").Value = ""
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H134").Value = 74475.57000000001
$ws.Range("J134").Value = 74475.57000000001
$ws.Range("L134").Value = 74475.57000000001
$ws.Range("N134").Value = -84615.57000000001
$ws.Range("H135").Value = 30153.5
$ws.Range("J135").Value = 30153.5
$ws.Range("L135").Value = 30153.5
$ws.Range("N135").Value = -40293.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 48852
$ws.Range("J46").Value = 48852
$ws.Range("L46").Value = 48852
$ws.Range("N46").Value = -49314
$ws.Range("H134").Value = 48852
$ws.Range("J134").Value = 48852
$ws.Range("L134").Value = 146556
$ws.Range("N134").Value = -151626
$ws.Range("H135").Value = 29475.545
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 29475.545
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 29475.545
$ws.Range("M135").Value = ""
$ws.Range("N135").Value = -39615.545
